$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the red-fill highlight placed on F7 (empty, formatted cell) - "quick data cleaning"
$ws.Range("F7").Clear()

# Delete row 18 (Active Living / COVID refund entry) entirely; row 19
# (Little Caesars, -15) shifts up to become the new row 18 - "added missing negatives"
$ws.Rows(18).Delete()

# Leave the selection where the user last touched the sheet (the row they deleted)
[void]$ws.Range("A18:XFD18").Select()
